$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 9093823
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 9093823
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 27281469
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -27281805

$ws.Range("H40").Value = 5003750
$ws.Range("I40").Value = 2657.9092
$ws.Range("K40").Value = 2657.9092
$ws.Range("M40").Value = -2482.9092

$ws.Range("H116").Value = 3020.7058
$ws.Range("J116").Value = 3040.7273
$ws.Range("L116").Value = 3040.7273
$ws.Range("N116").Value = -9924.7273

$ws.Range("H118").Value = 523.5
$ws.Range("I118").Value = 523.5
$ws.Range("K118").Value = 1570.5
$ws.Range("M118").Value = 86.5

$ws.Range("H132").Value = 2113.6316
$ws.Range("I132").Value = 2113.6316
$ws.Range("K132").Value = 6340.8948
$ws.Range("M132").Value = -3810.8948

$ws.Range("H137").Value = 2405.2173
$ws.Range("I137").Value = 2025.8572
$ws.Range("K137").Value = 6077.571599999999
$ws.Range("M137").Value = -3527.571599999999

$ws.Range("H138").Value = 3956.3823
$ws.Range("J138").Value = 3842.2856
$ws.Range("L138").Value = 11526.8568
$ws.Range("N138").Value = -21806.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6638.9644
$ws.Range("I32").Value = 6638.9644
$ws.Range("K32").Value = 6638.9644
$ws.Range("M32").Value = -6351.9644

$ws.Range("H37").Value = 4566.6665
$ws.Range("I37").Value = 100
$ws.Range("J37").Value = 13500
$ws.Range("K37").Value = 100
$ws.Range("L37").Value = 13500
$ws.Range("M37").Value = 173
$ws.Range("N37").Value = -14046

$ws.Range("H80").Value = 20000
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 20000
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H134").Value = 190000
$ws.Range("J134").Value = 190000
$ws.Range("L134").Value = 190000
$ws.Range("N134").Value = -200140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 298.2857
$ws.Range("I64").Value = 377.6
$ws.Range("K64").Value = 377.6
$ws.Range("M64").Value = -152.6

$ws.Range("H67").Value = 298.2857
$ws.Range("I67").Value = 377.6
$ws.Range("K67").Value = 377.6
$ws.Range("M67").Value = 402.4

$ws.Range("H99").Value = 3624.3704
$ws.Range("I99").Value = 2517.9
$ws.Range("K99").Value = 2517.9
$ws.Range("M99").Value = -1019.9

$ws.Range("H107").Value = 2881.0435
$ws.Range("I107").Value = 2722.75
$ws.Range("K107").Value = 2722.75
$ws.Range("M107").Value = -802.75

$ws.Range("H134").Value = 5153.696
$ws.Range("I134").Value = 4427.8
$ws.Range("J134").Value = 9993
$ws.Range("K134").Value = 13283.4
$ws.Range("L134").Value = 29979
$ws.Range("M134").Value = -10748.4
$ws.Range("N134").Value = -35049

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 83340750
$ws.Range("I31").Value = 111116000
$ws.Range("K31").Value = 111116000
$ws.Range("M31").Value = -111115705

$ws.Range("H34").Value = 83340750
$ws.Range("I34").Value = 111116000
$ws.Range("K34").Value = 111116000
$ws.Range("M34").Value = -111115798

$ws.Range("H99").Value = 7352.75
$ws.Range("I99").Value = 5926
$ws.Range("K99").Value = 5926
$ws.Range("M99").Value = -4428

$ws.Range("H125").Value = 90000
$ws.Range("J125").Value = 90000
$ws.Range("L125").Value = 90000
$ws.Range("N125").Value = -94920

$ws.Range("H126").Value = 7352.75
$ws.Range("I126").Value = 5926
$ws.Range("K126").Value = 17778
$ws.Range("M126").Value = -15308

$ws.Range("H132").Value = 5236.9375
$ws.Range("I132").Value = 3060.8462
$ws.Range("K132").Value = 9182.5386
$ws.Range("M132").Value = -6652.5386

$ws.Range("H134").Value = 4219.3887
$ws.Range("I134").Value = 3324.9
$ws.Range("J134").Value = 5337.5
$ws.Range("K134").Value = 9974.700000000001
$ws.Range("L134").Value = 16012.5
$ws.Range("M134").Value = -7439.700000000001
$ws.Range("N134").Value = -21082.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 19397292
$ws.Range("I4").Value = 18170130
$ws.Range("J4").Value = 23078778
$ws.Range("K4").Value = 54510390
$ws.Range("L4").Value = 69236334
$ws.Range("M4").Value = -54510278
$ws.Range("N4").Value = -69236558

$ws.Range("H14").Value = 213.6
$ws.Range("I14").Value = 213.6
$ws.Range("K14").Value = 640.8
$ws.Range("M14").Value = -467.8

$ws.Range("H86").Value = 405.375
$ws.Range("I86").Value = 290.66666
$ws.Range("K86").Value = 871.9999799999999
$ws.Range("M86").Value = 314.0000200000001

$ws.Range("H89").Value = 405.375
$ws.Range("I89").Value = 290.66666
$ws.Range("K89").Value = 2615.99994
$ws.Range("M89").Value = 3312.00006

$ws.Range("H104").Value = 8016.6665
$ws.Range("J104").Value = 8024
$ws.Range("L104").Value = 24072
$ws.Range("N104").Value = -29314

$ws.Range("H112").Value = 8663.5
$ws.Range("I112").Value = 10527
$ws.Range("K112").Value = 31581
$ws.Range("M112").Value = -30473

$ws.Range("H118").Value = 975
$ws.Range("I118").Value = 466.66666
$ws.Range("K118").Value = 1399.99998
$ws.Range("M118").Value = -156.9999800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 816.6923
$ws.Range("I16").Value = 591.9091
$ws.Range("J16").Value = 2053
$ws.Range("K16").Value = 591.9091
$ws.Range("L16").Value = 2053
$ws.Range("M16").Value = -421.9091
$ws.Range("N16").Value = -2393

$ws.Range("H61").Value = 72443.28999999999
$ws.Range("J61").Value = 836.3333
$ws.Range("L61").Value = 836.3333
$ws.Range("N61").Value = -1240.3333

$ws.Range("H100").Value = 6378.4
$ws.Range("I100").Value = 5084.5
$ws.Range("J100").Value = 7857.143
$ws.Range("K100").Value = 5084.5
$ws.Range("L100").Value = 7857.143
$ws.Range("M100").Value = -4543.5
$ws.Range("N100").Value = -8939.143

$ws.Range("H113").Value = 72443.28999999999
$ws.Range("J113").Value = 836.3333
$ws.Range("L113").Value = 836.3333
$ws.Range("N113").Value = -5176.3333

$ws.Range("H136").Value = 86959940
$ws.Range("I136").Value = 58826564
$ws.Range("K136").Value = 176479692
$ws.Range("M136").Value = -176477142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 35000
$ws.Range("I51").Value = 35000
$ws.Range("K51").Value = 35000
$ws.Range("M51").Value = -34490

$ws.Range("H52").Value = 30466.666
$ws.Range("J52").Value = 31400
$ws.Range("L52").Value = 31400
$ws.Range("N52").Value = -31852

$ws.Range("H58").Value = 35000
$ws.Range("I58").Value = 35000
$ws.Range("K58").Value = 35000
$ws.Range("M58").Value = -34692

$ws.Range("H70").Value = 35000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 35000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H113").Value = 1101.7693
$ws.Range("I113").Value = 862.7826
$ws.Range("K113").Value = 2588.3478
$ws.Range("M113").Value = -418.3478
